# Update NATMI ligand-receptor edge statistics with new TPM-derived values
# (Tnc-Itgb3 sheet: columns G-T recomputed after ligand/receptor expression update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 3.575558
$ws.Range("H2").Value = 10.726674
$ws.Range("I2").Value = 0.025194653521236
$ws.Range("J2").Value = 0.02519465352123599
$ws.Range("M2").Value = 2.598166333333333
$ws.Range("N2").Value = 7.794499
$ws.Range("O2").Value = 0.3466013321552429
$ws.Range("P2").Value = 0.3466013321552429
$ws.Range("Q2").Value = 9.289894418480666
$ws.Range("R2").Value = 83.609049766326
$ws.Range("S2").Value = 0.008732500473650179
$ws.Range("T2").Value = 0.008732500473650177

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 3.575558
$ws.Range("H3").Value = 10.726674
$ws.Range("I3").Value = 0.025194653521236
$ws.Range("J3").Value = 0.02519465352123599
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("O3").Value = 0.5780859172985858
$ws.Range("P3").Value = 0.5780859172985858
$ws.Range("Q3").Value = 15.49433495572666
$ws.Range("R3").Value = 139.44901460154
$ws.Range("S3").Value = 0.01456467439184376
$ws.Range("T3").Value = 0.01456467439184376

# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 3.575558
$ws.Range("H4").Value = 10.726674
$ws.Range("I4").Value = 0.025194653521236
$ws.Range("J4").Value = 0.02519465352123599
$ws.Range("M4").Value = 0.4692043333333333
$ws.Range("N4").Value = 1.407613
$ws.Range("O4").Value = 0.06259293136852516
$ws.Range("P4").Value = 0.06259293136852516
$ws.Range("Q4").Value = 1.677667307684666
$ws.Range("R4").Value = 15.099005769162
$ws.Range("S4").Value = 0.001577007218708495
$ws.Range("T4").Value = 0.001577007218708495

# Row 5: ECs -> Resolving-Mac
$ws.Range("G5").Value = 3.575558
$ws.Range("H5").Value = 10.726674
$ws.Range("I5").Value = 0.025194653521236
$ws.Range("J5").Value = 0.02519465352123599
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09534933333333333
$ws.Range("N5").Value = 0.286048
$ws.Range("O5").Value = 0.01271981917764605
$ws.Range("P5").Value = 0.01271981917764604
$ws.Range("Q5").Value = 0.3409270715946666
$ws.Range("R5").Value = 3.068343644351999
$ws.Range("S5").Value = 0.0003204714370335651
$ws.Range("T5").Value = 0.000320471437033565

# Row 6: FAPs -> ECs
$ws.Range("I6").Value = 0.7460690747908298
$ws.Range("J6").Value = 0.7460690747908298
$ws.Range("M6").Value = 2.598166333333333
$ws.Range("N6").Value = 7.794499
$ws.Range("O6").Value = 0.3466013321552429
$ws.Range("P6").Value = 0.3466013321552429
$ws.Range("Q6").Value = 275.0941951973448
$ws.Range("R6").Value = 2475.847756776103
$ws.Range("S6").Value = 0.2585885352023312
$ws.Range("T6").Value = 0.2585885352023312

# Row 7: FAPs -> FAPs
$ws.Range("I7").Value = 0.7460690747908298
$ws.Range("J7").Value = 0.7460690747908298
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("O7").Value = 0.5780859172985858
$ws.Range("P7").Value = 0.5780859172985858
$ws.Range("Q7").Value = 458.8213183870411
$ws.Range("R7").Value = 4129.39186548337
$ws.Range("S7").Value = 0.4312920254685641
$ws.Range("T7").Value = 0.4312920254685641

# Row 8: FAPs -> MuSCs
$ws.Range("I8").Value = 0.7460690747908298
$ws.Range("J8").Value = 0.7460690747908298
$ws.Range("M8").Value = 0.4692043333333333
$ws.Range("N8").Value = 1.407613
$ws.Range("O8").Value = 0.06259293136852516
$ws.Range("P8").Value = 0.06259293136852516
$ws.Range("Q8").Value = 49.67941690470678
$ws.Range("R8").Value = 447.114752142361
$ws.Range("S8").Value = 0.04669865039456148
$ws.Range("T8").Value = 0.04669865039456148

# Row 9: FAPs -> Resolving-Mac
$ws.Range("I9").Value = 0.7460690747908298
$ws.Range("J9").Value = 0.7460690747908298
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09534933333333333
$ws.Range("N9").Value = 0.286048
$ws.Range("O9").Value = 0.01271981917764605
$ws.Range("P9").Value = 0.01271981917764604
$ws.Range("Q9").Value = 10.09560003122844
$ws.Range("R9").Value = 90.86040028105599
$ws.Range("S9").Value = 0.009489863725373039
$ws.Range("T9").Value = 0.009489863725373039

# Row 10: MuSCs -> ECs
$ws.Range("G10").Value = 32.36130266666667
$ws.Range("H10").Value = 97.08390800000001
$ws.Range("I10").Value = 0.2280292497513723
$ws.Range("J10").Value = 0.2280292497513723
$ws.Range("M10").Value = 2.598166333333333
$ws.Range("N10").Value = 7.794499
$ws.Range("O10").Value = 0.3466013321552429
$ws.Range("P10").Value = 0.3466013321552429
$ws.Range("Q10").Value = 84.08004709134356
$ws.Range("R10").Value = 756.720423822092
$ws.Range("S10").Value = 0.07903524173418623
$ws.Range("T10").Value = 0.07903524173418623

# Row 11: MuSCs -> FAPs
$ws.Range("G11").Value = 32.36130266666667
$ws.Range("H11").Value = 97.08390800000001
$ws.Range("I11").Value = 0.2280292497513723
$ws.Range("J11").Value = 0.2280292497513723
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("O11").Value = 0.5780859172985858
$ws.Range("P11").Value = 0.5780859172985858
$ws.Range("Q11").Value = 140.2345768467422
$ws.Range("R11").Value = 1262.11119162068
$ws.Range("S11").Value = 0.1318204980134304
$ws.Range("T11").Value = 0.1318204980134304

# Row 12: MuSCs -> MuSCs
$ws.Range("G12").Value = 32.36130266666667
$ws.Range("H12").Value = 97.08390800000001
$ws.Range("I12").Value = 0.2280292497513723
$ws.Range("J12").Value = 0.2280292497513723
$ws.Range("M12").Value = 0.4692043333333333
$ws.Range("N12").Value = 1.407613
$ws.Range("O12").Value = 0.06259293136852516
$ws.Range("P12").Value = 0.06259293136852516
$ws.Range("Q12").Value = 15.18406344351156
$ws.Range("R12").Value = 136.656570991604
$ws.Range("S12").Value = 0.01427301917970393
$ws.Range("T12").Value = 0.01427301917970393

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("G13").Value = 32.36130266666667
$ws.Range("H13").Value = 97.08390800000001
$ws.Range("I13").Value = 0.2280292497513723
$ws.Range("J13").Value = 0.2280292497513723
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09534933333333333
$ws.Range("N13").Value = 0.286048
$ws.Range("O13").Value = 0.01271981917764605
$ws.Range("P13").Value = 0.01271981917764604
$ws.Range("Q13").Value = 3.085628635064889
$ws.Range("R13").Value = 27.770657715584
$ws.Range("S13").Value = 0.002900490824051745
$ws.Range("T13").Value = 0.002900490824051745

# Row 14: Resolving-Mac -> ECs
$ws.Range("G14").Value = 0.1003386666666667
$ws.Range("H14").Value = 0.301016
$ws.Range("I14").Value = 0.000707021936561918
$ws.Range("J14").Value = 0.0007070219365619179
$ws.Range("M14").Value = 2.598166333333333
$ws.Range("N14").Value = 7.794499
$ws.Range("O14").Value = 0.3466013321552429
$ws.Range("P14").Value = 0.3466013321552429
$ws.Range("Q14").Value = 0.2606965456648889
$ws.Range("R14").Value = 2.346268910984
$ws.Range("S14").Value = 0.0002450547450753404
$ws.Range("T14").Value = 0.0002450547450753404

# Row 15: Resolving-Mac -> FAPs
$ws.Range("G15").Value = 0.1003386666666667
$ws.Range("H15").Value = 0.301016
$ws.Range("I15").Value = 0.000707021936561918
$ws.Range("J15").Value = 0.0007070219365619179
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("O15").Value = 0.5780859172985858
$ws.Range("P15").Value = 0.5780859172985858
$ws.Range("Q15").Value = 0.4348079125955556
$ws.Range("R15").Value = 3.91327121336
$ws.Range("S15").Value = 0.000408719424747619
$ws.Range("T15").Value = 0.0004087194247476189

# Row 16: Resolving-Mac -> MuSCs
$ws.Range("G16").Value = 0.1003386666666667
$ws.Range("H16").Value = 0.301016
$ws.Range("I16").Value = 0.000707021936561918
$ws.Range("J16").Value = 0.0007070219365619179
$ws.Range("M16").Value = 0.4692043333333333
$ws.Range("N16").Value = 1.407613
$ws.Range("O16").Value = 0.06259293136852516
$ws.Range("P16").Value = 0.06259293136852516
$ws.Range("Q16").Value = 0.04707933720088889
$ws.Range("R16").Value = 0.423714034808
$ws.Range("S16").Value = 0.00004425457555126188
$ws.Range("T16").Value = 0.00004425457555126188

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Range("G17").Value = 0.1003386666666667
$ws.Range("H17").Value = 0.301016
$ws.Range("I17").Value = 0.000707021936561918
$ws.Range("J17").Value = 0.0007070219365619179
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09534933333333333
$ws.Range("N17").Value = 0.286048
$ws.Range("O17").Value = 0.01271981917764605
$ws.Range("P17").Value = 0.01271981917764604
$ws.Range("Q17").Value = 0.009567224974222223
$ws.Range("R17").Value = 0.086105024768
$ws.Range("S17").Value = 0.00000899319118769673
$ws.Range("T17").Value = 0.000008993191187696729
